# Scale the "Commerce intra-continental" (O) and "Commerce extra-continental" (P)
# columns from billions to millions of dollars (x1000) for data rows 3-99 on Tab21.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab21")

for ($r = 3; $r -le 99; $r++) {
    $oCell = $ws.Cells.Item($r, 15)  # column O
    $pCell = $ws.Cells.Item($r, 16)  # column P

    if ($oCell.Value2 -ne $null) {
        $oCell.Value = $oCell.Value2 * 1000
    }
    if ($pCell.Value2 -ne $null) {
        $pCell.Value = $pCell.Value2 * 1000
    }
}

# Bump the saved window height (cosmetic view-state change recorded in workbook.xml)
$excel.ActiveWindow.Height = 12490 / 20
